# Updated symbol list on Tue Jan 17 18:14:47 UTC 2023 with GitHub Actions
#
# Refreshes the crypto snapshot table on the active sheet:
#   - Column D (Price) and column E (Volume(1h)) get new quoted readings
#     for the rows whose numbers moved since the last run.
#   - Column G (Hora) advances from "17" to "18" for every data row
#     (rows 2-51), since the sheet was re-scraped an hour later.
#
# Values are written with a leading "'" (quote prefix) so Excel stores
# them as literal text (matching the sheet's existing text-formatted
# Price/Volume/Hora columns) instead of re-interpreting numeric- or
# percent-looking strings as actual numbers/percentages.

$updates = @{
    2  = @{ D = "301.92";     E = "1.08%" }
    3  = @{ D = "31.98";      E = "0.93%" }
    4  = @{ D = "5.034";      E = "-0.92%" }
    5  = @{ D = "0.07822";    E = "-4.03%" }
    6  = @{ D = "2.093";      E = "-18.78%" }
    7  = @{ D = "7.800";      E = "0.37%" }
    8  = @{ D = "3.771";      E = "-1.82%" }
    9  = @{ D = "0.9243";     E = "-0.59%" }
    10 = @{ D = "0.1748";     E = "-0.82%" }
    11 = @{ D = "0.07890";    E = "5.27%" }
    12 = @{ D = "0.08781";    E = "-3.92%" }
    13 = @{ D = "0.03126";    E = "3.22%" }
    14 = @{ D = "0.09999" }
    15 = @{ D = "0.001511";   E = "1.24%" }
    16 = @{ D = "0.005913";   E = "-0.01%" }
    17 = @{ D = "3.450";      E = "-3.44%" }
    18 = @{ D = "2.269";      E = "0.97%" }
    19 = @{ E = "1.37%" }
    20 = @{ E = "-1.36%" }
    21 = @{ D = "4.154";      E = "4.93%" }
    22 = @{ D = "0.1791";     E = "5.46%" }
    23 = @{ D = "0.04591";    E = "-0.30%" }
    24 = @{ E = "-0.17%" }
    25 = @{ D = "0.004481";   E = "0.31%" }
    26 = @{ D = "0.0001249";  E = "4.20%" }
    39 = @{ D = "0.01730";    E = "-2.15%" }
    40 = @{ D = "0.04778";    E = "5.24%" }
    41 = @{ D = "0.007389";   E = "7.21%" }
    42 = @{ E = "0.27%" }
    43 = @{ D = "0.002078";   E = "-5.85%" }
    44 = @{ D = "0.01073";    E = "8.08%" }
    45 = @{ D = "0.00006079"; E = "-5.77%" }
    46 = @{ D = "0.00000000749"; E = "0.05%" }
    47 = @{ D = "0.003497";   E = "-59.98%" }
    49 = @{ D = "0.00002098"; E = "0.05%" }
    50 = @{ D = "0.0001998";  E = "0.05%" }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $ws.Cells.Item($row, 4).Value = "'" + $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Cells.Item($row, 5).Value = "'" + $vals["E"]
    }
}

# Every data row's "Hora" column moves from 17 -> 18.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 7).Value = "'18"
}
